# Updated cryptos list on Wed May 10 19:54:58 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns hold plain text (inline strings) in the source
# file even though most of the values look numeric (e.g. "1.001"). Force
# the whole data range to Text format first so Excel's auto-detection
# doesn't silently convert these assignments into real numbers; the
# formatting is reverted to the sheet's original (unstyled) look at the
# very end so no extra cell styling is introduced.
$ws.Range("D2:E51").NumberFormat = "@"

# Rows 46/47 swap position (Decentraland moves up to #46, Quant moves
# down to #47) and both carry refreshed price/volume figures.
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.4792"
$ws.Range("E46").Value = "  +2.87%  "

$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "106.98"
$ws.Range("E47").Value = "  +2.27%  "

# Refreshed price (D) / volume-1h (E) figures for the rest of the rows.
$ws.Range("D2").Value = "27.723.19"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.851.56"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "315.37"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "0.4311"
$ws.Range("E7").Value = "  +1.36%  "
$ws.Range("D8").Value = "0.3709"
$ws.Range("E8").Value = "  +1.78%  "
$ws.Range("D9").Value = "0.07361"
$ws.Range("E9").Value = "  +0.86%  "
$ws.Range("D10").Value = "0.8751"
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").Value = "21.01"
$ws.Range("E11").Value = "  +2.06%  "
$ws.Range("D12").Value = "1.883.16"
$ws.Range("E12").Value = "  +4.05%  "
$ws.Range("D13").Value = "5.451"
$ws.Range("E13").Value = "  +2.63%  "
$ws.Range("D14").Value = "6.598"
$ws.Range("E14").Value = "  +1.27%  "
$ws.Range("D15").Value = "0.06942"
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "81.24"
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("D18").Value = "0.000009089"
$ws.Range("E18").Value = "  +1.07%  "
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "15.58"
$ws.Range("E20").Value = "  +1.53%  "
$ws.Range("D21").Value = "27.753.18"
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("D22").Value = "5.091"
$ws.Range("E22").Value = "  +2.22%  "
$ws.Range("D23").Value = "11.01"
$ws.Range("E23").Value = "  +5.98%  "
$ws.Range("D24").Value = "2.093.01"
$ws.Range("E24").Value = "  +2.75%  "
$ws.Range("D25").Value = "1.966"
$ws.Range("E25").Value = "  -0.66%  "
$ws.Range("D26").Value = "155.05"
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("D27").Value = "18.59"
$ws.Range("E27").Value = "  -1.42%  "
$ws.Range("D28").Value = "5.320"
$ws.Range("E28").Value = "  +1.23%  "
$ws.Range("D29").Value = "115.44"
$ws.Range("E29").Value = "  -5.32%  "
$ws.Range("D30").Value = "1.851"
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("D31").Value = "0.08930"
$ws.Range("D32").Value = "0.7828"
$ws.Range("E32").Value = "  +1.85%  "
$ws.Range("D33").Value = "4.606"
$ws.Range("E33").Value = "  +1.60%  "
$ws.Range("D34").Value = "2.974"
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("E35").Value = "  +5.06%  "
$ws.Range("D36").Value = "1.000"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("E37").Value = "  +1.77%  "
$ws.Range("D38").Value = "0.05426"
$ws.Range("E38").Value = "  +0.81%  "
$ws.Range("D39").Value = "0.01960"
$ws.Range("E39").Value = "  +0.97%  "
$ws.Range("D40").Value = "2.840"
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("D41").Value = "0.5214"
$ws.Range("E41").Value = "  +3.02%  "
$ws.Range("E42").Value = "  +2.11%  "
$ws.Range("D43").Value = "6.763"
$ws.Range("E43").Value = "  -0.81%  "
$ws.Range("D44").Value = "8.666"
$ws.Range("E44").Value = "  +3.47%  "
$ws.Range("D45").Value = "10.68"
$ws.Range("E45").Value = "  +3.06%  "
$ws.Range("D48").Value = "0.06564"
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("D49").Value = "0.9999"
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").Value = "1.670"
$ws.Range("E50").Value = "  +2.67%  "
$ws.Range("D51").Value = "1.841"
$ws.Range("E51").Value = "  +5.69%  "

# Restore the original (unstyled / default) look of the data range now
# that every value is safely stored as text.
$ws.Range("D2:E51").Style = "Normal"
